$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at the top for a "100 Iterations" title, pushing
#    the existing table down one row (old row N -> new row N+1).
# ---------------------------------------------------------------------
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "100 Iterations"

# ---------------------------------------------------------------------
# 2. Update the (now shifted) first table with the final run's numbers.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5

$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5

$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 3

# Replace the static H/I/J values with real formulas referencing C/D/E.
$ws.Range("H4").Formula = "=C4/10"
$ws.Range("I4").Formula = "=D4/10"
$ws.Range("J4").Formula = "=E4/10"

$ws.Range("H5").Formula = "=C5/10"
$ws.Range("I5").Formula = "=D5/10"
$ws.Range("J5").Formula = "=E5/10"

$ws.Range("H6").Formula = "=C6/9"
$ws.Range("I6").Formula = "=D6/9"
$ws.Range("J6").Formula = "=E6/9"

# The bottom data row's H/I/J keeps a 2-decimal display format.
$ws.Range("H6:J6").NumberFormat = "0.00"

# Clear the stray "applyNumberFormat"-only style that had no visible
# effect on the old C5 (now C6) cell - it's no longer applied anywhere.
$ws.Range("C6").ClearFormats()

# ---------------------------------------------------------------------
# 3. Build a second, identical-structure table (rows 8-14) with the
#    second classification run's numbers ("100 Iterations" again).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "100 Iterations"

$ws.Range("B9").Value = "Control"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("C9:E9").HorizontalAlignment = -4108

$ws.Range("A10").Value = "Atypical"
$ws.Range("C10:E10").HorizontalAlignment = -4108

$ws.Range("C9:C10").Merge()
$ws.Range("D9:D10").Merge()
$ws.Range("E9:E10").Merge()

$ws.Range("A11").Value = 0
$ws.Range("A11:B11").HorizontalAlignment = -4108
$ws.Range("A11:B11").Merge()
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5
$ws.Range("F11").Formula = "=SUM(C11:E11)/30"
$ws.Range("H11").Formula = "=C11/10"
$ws.Range("I11").Formula = "=D11/10"
$ws.Range("J11").Formula = "=E11/10"

$ws.Range("A12").Value = 1
$ws.Range("A12:B12").HorizontalAlignment = -4108
$ws.Range("A12:B12").Merge()
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 3
$ws.Range("F12").Formula = "=SUM(C12:E12)/30"
$ws.Range("H12").Formula = "=C12/10"
$ws.Range("I12").Formula = "=D12/10"
$ws.Range("J12").Formula = "=E12/10"

$ws.Range("A13").Value = 2
$ws.Range("A13:B13").HorizontalAlignment = -4108
$ws.Range("A13:B13").Merge()
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 2
$ws.Range("F13").Formula = "=SUM(C13:E13)/27"
$ws.Range("H13").Formula = "=C13/9"
$ws.Range("I13").Formula = "=D13/9"
$ws.Range("J13").Formula = "=E13/9"

$ws.Range("C14").Formula = "=SUM(C11:C13)/29"
$ws.Range("D14").Formula = "=SUM(D11:D13)/29"
$ws.Range("E14").Formula = "=SUM(E11:E13)/29"
$ws.Range("F14").Formula = "=SUM(C11:E13)/87"

# ---------------------------------------------------------------------
# 4. Move the connector/line drawing down by one row, matching the new
#    title row inserted at the top (its anchor moves from row 0 to 1,
#    and from row 2 to row 3).
# ---------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + $ws.Rows(1).RowHeight

# ---------------------------------------------------------------------
# 5. Final view bits: reflect the new selection and used range.
# ---------------------------------------------------------------------
$ws.Range("G14").Select()

Write-Host "done"
